$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$sub5 = [char]0x2085

# Row 2
Set-CellText $ws "D2" "26.024.99"
Set-CellText $ws "E2" "  -1.87%  "

# Row 3
Set-CellText $ws "D3" "1.639.02"
Set-CellText $ws "E3" "  -1.92%  "

# Row 4
Set-CellText $ws "E4" "  -0.31%  "

# Row 5
Set-CellText $ws "D5" "215.36"
Set-CellText $ws "E5" "  -2.11%  "

# Row 6
Set-CellText $ws "D6" "0.5048"
Set-CellText $ws "E6" "  -2.53%  "

# Row 7
Set-CellText $ws "D7" "1.008"
Set-CellText $ws "E7" "  -0.16%  "

# Row 8
Set-CellText $ws "D8" "0.2576"
Set-CellText $ws "E8" "  +0.17%  "

# Row 9
Set-CellText $ws "D9" "0.06395"
Set-CellText $ws "E9" "  -1.35%  "

# Row 10
Set-CellText $ws "D10" "19.53"
Set-CellText $ws "E10" "  -1.95%  "

# Row 11
Set-CellText $ws "D11" "0.07725"
Set-CellText $ws "E11" "  +0.42%  "

# Row 12
Set-CellText $ws "D12" "1.648.41"
Set-CellText $ws "E12" "  -1.67%  "

# Row 13
Set-CellText $ws "D13" "4.244"
Set-CellText $ws "E13" "  -2.09%  "

# Row 14
Set-CellText $ws "D14" "1.867.08"
Set-CellText $ws "E14" "  -2.13%  "

# Row 15
Set-CellText $ws "D15" "0.5449"
Set-CellText $ws "E15" "  -2.56%  "

# Row 16
Set-CellText $ws "D16" "0.0${sub5}7937"
Set-CellText $ws "E16" "  -0.93%  "

# Row 17
Set-CellText $ws "D17" "63.61"
Set-CellText $ws "E17" "  -2.56%  "

# Row 18
Set-CellText $ws "D18" "26.041.76"
Set-CellText $ws "E18" "  -2.20%  "

# Row 19
Set-CellText $ws "E19" "  -0.24%  "

# Row 20
Set-CellText $ws "D20" "205.66"
Set-CellText $ws "E20" "  -3.14%  "

# Row 21
Set-CellText $ws "D21" "4.341"
Set-CellText $ws "E21" "  -3.15%  "

# Row 22
Set-CellText $ws "D22" "9.978"
Set-CellText $ws "E22" "  -1.40%  "

# Row 23
Set-CellText $ws "D23" "5.981"
Set-CellText $ws "E23" "  +1.20%  "

# Row 24
Set-CellText $ws "E24" "  -0.25%  "

# Row 25
Set-CellText $ws "D25" "1.948"
Set-CellText $ws "E25" "  +12.65%  "

# Row 26
Set-CellText $ws "D26" "142.29"
Set-CellText $ws "E26" "  -0.44%  "

# Row 27
Set-CellText $ws "D27" "0.1157"
Set-CellText $ws "E27" "  -0.86%  "

# Row 28
Set-CellText $ws "D28" "6.861"
Set-CellText $ws "E28" "  -2.28%  "

# Row 29
Set-CellText $ws "D29" "15.77"
Set-CellText $ws "E29" "  -0.05%  "

# Row 30
Set-CellText $ws "B30" "Hedera"
Set-CellText $ws "C30" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText $ws "D30" "0.05005"
Set-CellText $ws "E30" "  -4.60%  "

# Row 31
Set-CellText $ws "B31" "PancakeSwap"
Set-CellText $ws "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText $ws "D31" "1.238"
Set-CellText $ws "E31" "  -2.54%  "

# Row 32
Set-CellText $ws "D32" "3.274"
Set-CellText $ws "E32" "  -2.43%  "

# Row 33
Set-CellText $ws "D33" "3.199"
Set-CellText $ws "E33" "  -1.25%  "

# Row 34
Set-CellText $ws "D34" "1.536"
Set-CellText $ws "E34" "  -3.70%  "

# Row 35
Set-CellText $ws "D35" "2.334"
Set-CellText $ws "E35" "  -2.45%  "

# Row 36
Set-CellText $ws "D36" "0.9066"
Set-CellText $ws "E36" "  -2.43%  "

# Row 37
Set-CellText $ws "D37" "2.650"
Set-CellText $ws "E37" "  -4.40%  "

# Row 38
Set-CellText $ws "D38" "0.5669"
Set-CellText $ws "E38" "  -1.04%  "

# Row 39
Set-CellText $ws "D39" "1.123.88"
Set-CellText $ws "E39" "  -3.52%  "

# Row 40
Set-CellText $ws "D40" "0.01561"
Set-CellText $ws "E40" "  -2.61%  "

# Row 41
Set-CellText $ws "E41" "  -0.16%  "

# Row 42
Set-CellText $ws "D42" "2.553"
Set-CellText $ws "E42" "  -1.17%  "

# Row 43
Set-CellText $ws "D43" "5.633"
Set-CellText $ws "E43" "  -1.33%  "

# Row 44
Set-CellText $ws "D44" "0.8116"
Set-CellText $ws "E44" "  -2.42%  "

# Row 45
Set-CellText $ws "D45" "99.66"
Set-CellText $ws "E45" "  +0.04%  "

# Row 46
Set-CellText $ws "D46" "1.778.45"
Set-CellText $ws "E46" "  -2.26%  "

# Row 47
Set-CellText $ws "E47" "  -0.46%  "

# Row 48
Set-CellText $ws "D48" "0.4532"
Set-CellText $ws "E48" "  +0.79%  "

# Row 49
Set-CellText $ws "E49" "  +0.13%  "

# Row 50
Set-CellText $ws "D50" "54.87"
Set-CellText $ws "E50" "  -1.97%  "

# Row 51
Set-CellText $ws "D51" "7.730"
Set-CellText $ws "E51" "  -3.07%  "
